$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.139.49"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.41%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.677.45"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.07%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "214.06"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.83%  "
$ws.Range("E6").Value = "  +0.08%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "22.77"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +6.45%  "
$ws.Range("E9").Value = "  +2.05%  "
$ws.Range("E10").Value = "  -0.70%  "
$ws.Range("E11").Value = "  +0.23%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.915.26"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.688.09"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.35%  "
$ws.Range("E14").Value = "  +2.26%  "
$ws.Range("E15").Value = "  +2.79%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "66.48"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.03%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "27.098.85"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.29%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "234.87"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.88"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.38%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0740"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.22%  "
$ws.Range("E21").Value = "  +0.08%  "
$ws.Range("E22").Value = "  +1.39%  "
$ws.Range("E23").Value = "  +2.54%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.08"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.66%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "147.97"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.92%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.43"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.23%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.32"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.56%  "
$ws.Range("E28").Value = "  -0.07%  "
$ws.Range("E29").Value = "  -0.02%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0500"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.48%  "
$ws.Range("E31").Value = "  -0.46%  "
$ws.Range("E32").Value = "  -0.34%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.542.52"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.22%  "
$ws.Range("E34").Value = "  +1.29%  "
$ws.Range("E35").Value = "  -3.29%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.604"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.53%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.938"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.03%  "
$ws.Range("E38").Value = "  -0.49%  "
$ws.Range("E39").Value = "  -1.55%  "
$ws.Range("E40").Value = "  +2.55%  "
$ws.Range("E41").Value = "  +3.48%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "69.28"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.88%  "
$ws.Range("E43").Value = "  +0.06%  "
$ws.Range("E44").Value = "  -0.11%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.821.92"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.778"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.29%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "89.84"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.65%  "
$ws.Range("E48").Value = "  +5.76%  "
$ws.Range("E49").Value = "  +2.22%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.23"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.29%  "
$ws.Range("E51").Value = "  -0.48%  "
